# Applies the "Automatic update of files." commit:
#   1) Bumps the "Förändrad" date (column C) from 2023-09-15 (45184) to
#      2023-09-17 (45186) for every data row.
#   2) Adds a second argument (the case id, taken from column A) to every
#      HYPERLINK() formula in columns S-Y so the link shows a friendly
#      label instead of the raw URL.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 292
$oldDate = 45184
$newDate = 45186

# --- 1) Column C: refresh the "changed" date stamp on every data row ----
for ($r = 2; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $current = $cCell.Value2
    if ($current -eq "$oldDate") {
        $cCell.Value = $newDate
    }
}

# --- 2) Columns S-Y: rebuild HYPERLINK formulas with a label argument ---
# col index -> (subfolder, extension) beneath the GitHub Pages base URL
$baseUrl = "https://klasma.github.io/Logging_VALDEMARSVIK/"
$cols      = @(19,        20,       21,        22,          23,                24,         25)
$subfolder = @("artfynd", "kartor", "knärot",  "klagomål",  "klagomålsmail",   "tillsyn",  "tillsynsmail")
$extension = @("xlsx",    "png",    "png",     "docx",      "docx",            "docx",     "docx")

for ($r = 2; $r -le $lastRow; $r++) {
    $idCell = $ws.Cells.Item($r, 1)
    $caseId = $idCell.Text
    if ([string]::IsNullOrEmpty($caseId)) { continue }

    for ($k = 0; $k -lt $cols.Length; $k++) {
        $col = $cols[$k]
        $cell = $ws.Cells.Item($r, $col)
        $existingFormula = $cell.Formula
        if ([string]::IsNullOrEmpty($existingFormula)) { continue }

        $url = $baseUrl + $subfolder[$k] + "/" + $caseId + "." + $extension[$k]
        $newFormula = "=HYPERLINK(""" + $url + """, """ + $caseId + """)"
        $cell.Formula = $newFormula
    }
}
